$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '50.147.81'
$ws.Range('E2').Value = '  +3.82%  '

# Row 3
$ws.Range('D3').Value = '2.608.50'
$ws.Range('E3').Value = '  +3.79%  '

# Row 4
$ws.Range('D4').Value = '''0.997'
$ws.Range('E4').Value = '  -0.28%  '

# Row 5
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').Value = '''110.53'
$ws.Range('E5').Value = '  +1.56%  '

# Row 6
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '''324.10'
$ws.Range('E6').Value = '  +0.80%  '

# Row 7
$ws.Range('D7').Value = '''0.533'
$ws.Range('E7').Value = '  +0.95%  '

# Row 8
$ws.Range('E8').Value = '  -0.19%  '

# Row 9
$ws.Range('D9').Value = '''0.563'
$ws.Range('E9').Value = '  +3.35%  '

# Row 10
$ws.Range('D10').Value = '''40.83'
$ws.Range('E10').Value = '  +1.97%  '

# Row 11
$ws.Range('D11').Value = '''20.80'
$ws.Range('E11').Value = '  +3.81%  '

# Row 12
$ws.Range('D12').Value = '''0.0825'
$ws.Range('E12').Value = '  +0.68%  '

# Row 13
$ws.Range('E13').Value = '  +0.59%  '

# Row 14
$ws.Range('D14').Value = '''7.31'
$ws.Range('E14').Value = '  +1.44%  '

# Row 15
$ws.Range('D15').Value = '3.015.56'
$ws.Range('E15').Value = '  +3.63%  '

# Row 16
$ws.Range('D16').Value = '2.573.84'
$ws.Range('E16').Value = '  +2.46%  '

# Row 17
$ws.Range('D17').Value = '''0.869'
$ws.Range('E17').Value = '  +2.77%  '

# Row 18
$ws.Range('D18').Value = '49.947.29'
$ws.Range('E18').Value = '  +3.76%  '

# Row 19
$ws.Range('D19').Value = '''3.08'
$ws.Range('E19').Value = '  +12.05%  '

# Row 20
$ws.Range('D20').Value = '''13.46'
$ws.Range('E20').Value = '  +2.33%  '

# Row 21
$ws.Range('D21').Value = '''6.78'
$ws.Range('E21').Value = '  +0.16%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0956'
$ws.Range('E22').Value = '  +0.58%  '

# Row 23
$ws.Range('D23').Value = '''285.80'
$ws.Range('E23').Value = '  +2.75%  '

# Row 24
$ws.Range('D24').Value = '''73.12'
$ws.Range('E24').Value = '  +1.15%  '

# Row 25
$ws.Range('D25').Value = '''2.55'
$ws.Range('E25').Value = '  -0.59%  '

# Row 26
$ws.Range('D26').Value = '''26.81'
$ws.Range('E26').Value = '  +3.43%  '

# Row 27
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.03%  '

# Row 28
$ws.Range('D28').Value = '''0.147'
$ws.Range('E28').Value = '  +5.28%  '

# Row 29
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '''10.01'
$ws.Range('E29').Value = '  +1.42%  '

# Row 30
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''2.23'
$ws.Range('E30').Value = '  -7.07%  '

# Row 31
$ws.Range('D31').Value = '''36.12'
$ws.Range('E31').Value = '  +1.84%  '

# Row 32
$ws.Range('D32').Value = '''49.56'
$ws.Range('E32').Value = '  +0.75%  '

# Row 33
$ws.Range('D33').Value = '''19.93'
$ws.Range('E33').Value = '  +2.38%  '

# Row 34
$ws.Range('D34').Value = '''5.45'
$ws.Range('E34').Value = '  +1.43%  '

# Row 35
$ws.Range('E35').Value = '  -0.44%  '

# Row 36
$ws.Range('D36').Value = '''0.0794'
$ws.Range('E36').Value = '  +1.02%  '

# Row 37
$ws.Range('D37').Value = '''2.07'
$ws.Range('E37').Value = '  +5.37%  '

# Row 38
$ws.Range('D38').Value = '''4.77'
$ws.Range('E38').Value = '  +2.74%  '

# Row 39
$ws.Range('D39').Value = '''3.07'
$ws.Range('E39').Value = '  +3.82%  '

# Row 40
$ws.Range('D40').Value = '''124.31'
$ws.Range('E40').Value = '  +1.54%  '

# Row 41
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '''22.75'
$ws.Range('E41').Value = '  +5.02%  '

# Row 42
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '''0.113'
$ws.Range('E42').Value = '  +0.73%  '

# Row 43
$ws.Range('E43').Value = '  +0.26%  '

# Row 44
$ws.Range('D44').Value = '''0.0315'
$ws.Range('E44').Value = '  +2.78%  '

# Row 45
$ws.Range('D45').Value = '''3.38'
$ws.Range('E45').Value = '  +6.73%  '

# Row 46
$ws.Range('D46').Value = '2.040.43'
$ws.Range('E46').Value = '  +2.02%  '

# Row 47
$ws.Range('D47').Value = '''2.05'
$ws.Range('E47').Value = '  +10.57%  '

# Row 48
$ws.Range('D48').Value = '''2.16'
$ws.Range('E48').Value = '  +8.88%  '

# Row 49
$ws.Range('D49').Value = '''9.20'
$ws.Range('E49').Value = '  +1.83%  '

# Row 50
$ws.Range('D50').Value = '''5.42'
$ws.Range('E50').Value = '  +3.49%  '

# Row 51
$ws.Range('D51').Value = '''81.91'
$ws.Range('E51').Value = '  +2.05%  '
